$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (data rows 2-51) to Text format so numeric-looking
# strings (e.g. "28.678.45", "0.07510") are preserved exactly as text,
# matching the inline-string cell type used in the source workbook.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.678.45"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.805.38"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "317.60"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "0.5424"
$ws.Range("E7").Value = "  -2.68%  "

$ws.Range("D8").Value = "0.3798"
$ws.Range("E8").Value = "  -1.43%  "

$ws.Range("D9").Value = "0.07510"
$ws.Range("E9").Value = "  -0.99%  "

$ws.Range("D10").Value = "42.42"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("D11").Value = "1.113"
$ws.Range("E11").Value = "  -1.45%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "20.63"
$ws.Range("E13").Value = "  -2.16%  "

$ws.Range("D14").Value = "6.157"
$ws.Range("E14").Value = "  -1.40%  "

$ws.Range("D15").Value = "7.363"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "1.802.59"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").Value = "90.13"
$ws.Range("E17").Value = "  -1.22%  "

$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "0.06496"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").Value = "17.35"
$ws.Range("E21").Value = "  +0.59%  "

$ws.Range("D22").Value = "5.953"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("D23").Value = "28.696.78"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").Value = "11.14"
$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("D25").Value = "2.098"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("D26").Value = "161.66"
$ws.Range("E26").Value = "  +2.60%  "

$ws.Range("D27").Value = "20.48"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").Value = "2.007.49"
$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("D29").Value = "2.337"
$ws.Range("E29").Value = "  -4.10%  "

$ws.Range("D30").Value = "122.88"
$ws.Range("E30").Value = "  -0.90%  "

$ws.Range("D31").Value = "1.138"
$ws.Range("E31").Value = "  -2.04%  "

$ws.Range("D32").Value = "0.1056"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.647"
$ws.Range("E33").Value = "  -1.95%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.687"
$ws.Range("E34").Value = "  +1.13%  "

$ws.Range("D35").Value = "0.06668"
$ws.Range("E35").Value = "  +7.42%  "

$ws.Range("D36").Value = "0.2259"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").Value = "0.02310"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").Value = "8.633"
$ws.Range("E38").Value = "  -2.80%  "

$ws.Range("D39").Value = "5.024"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("E40").Value = "  -3.10%  "

$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("D42").Value = "1.198"
$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("D43").Value = "1.450"
$ws.Range("E43").Value = "  +4.97%  "

$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "13.30"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "3.703"
$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").Value = "0.5849"
$ws.Range("E47").Value = "  -2.40%  "

$ws.Range("D48").Value = "126.76"
$ws.Range("E48").Value = "  +2.86%  "

$ws.Range("D49").Value = "1.952"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").Value = "1.160"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("D51").Value = "0.06910"
$ws.Range("E51").Value = "  -0.37%  "
